$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.717.43"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "2.677.73"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.88"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.76"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("E8").Value = "  +5.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("E14").Value = "  -2.76%  "

$ws.Range("D15").Value = "3.157.53"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "65.555.37"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "2.681.95"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.91"
$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.48"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.73"
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +5.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").Value = "  -1.58%  "

$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("E28").Value = "  -6.17%  "

$ws.Range("E29").Value = "  -2.77%  "

$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.34"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.14"
$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.58"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.45"
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("E40").Value = "  -3.13%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "164.78"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.33"
$ws.Range("E44").Value = "  +2.18%  "

$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.91"
$ws.Range("E46").Value = "  -3.23%  "

$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("E48").Value = "  -3.18%  "

$ws.Range("D49").Value = "0.0₆0260"
$ws.Range("E49").Value = "  +13.51%  "

$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.10"
$ws.Range("E51").Value = "  -4.77%  "
